# Q3 Update - 2025
# Applies the shared-string / cell-value corrections described by the diff:
#  - the "short-url" column (B) text changes for every data row
#  - a handful of numeric-looking "asylum_seekers"/"refugees"/"oip" values
#    get corrected for rows 544-587

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Every data row shares the same "short-url" shared string (B2:B587).
$ws.Range("B2:B587").Value = "XKn14k"

# 2) Numeric-looking text cells need to stay text (shared string), not be
#    re-interpreted as numbers, so format each target cell as Text before
#    writing the new value.
$cellUpdates = @{
    "O544" = "11"
    "O545" = "38"
    "O546" = "49"
    "O547" = "16"
    "N548" = "6"
    "O548" = "56"
    "O549" = "20"
    "N552" = "175"
    "O552" = "292"
    "N553" = "47"
    "O553" = "23825"
    "O554" = "1407"
    "O555" = "68"
    "O560" = "24"
    "O570" = "42"
    "O572" = "9"
    "O573" = "103"
    "O575" = "24"
    "O576" = "12"
    "O577" = "208"
    "O581" = "9"
    "O584" = "6"
    "V585" = "7791"
    "O586" = "5"
    "N587" = "799"
    "O587" = "4416"
    "U587" = "35951"
}

foreach ($ref in $cellUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $cellUpdates[$ref]
}
